# Update Active_Outages.xlsx - 6/18/2025, 4:35:15 PM
$wb = $excel.ActiveWorkbook

# --- Sheet R1: refresh elapsed-duration figures for existing rows ---
$wsR1 = $wb.Worksheets.Item("R1")
$wsR1.Range("G2").Value = "3929:49:11"
$wsR1.Range("G3").Value = "69:21:49"

# --- Sheet R1, row 4: new outage record filled in ---
$wsR1.Range("C4").Value = "1"
$wsR1.Range("D4").Value = "JED0925"
$wsR1.Range("E4").Value = "Critical"
$wsR1.Range("F4").Value = "2025-06-14 20:13:06"
$wsR1.Range("G4").Value = "93:21:49"
$wsR1.Range("J4").Value = "In progress"
$wsR1.Range("K4").Value = "team on the way"

# --- Sheet R2: refresh elapsed-duration figures ---
$wsR2 = $wb.Worksheets.Item("R2")
$wsR2.Range("G2").Value = "12111:12:24"
$wsR2.Range("G3").Value = "3240:55:53"
$wsR2.Range("G4").Value = "479:07:27"

# --- Sheet R4: refresh elapsed-duration figures ---
$wsR4 = $wb.Worksheets.Item("R4")
$wsR4.Range("G2").Value = "2957:02:13"
$wsR4.Range("G3").Value = "184:14:28"
$wsR4.Range("G4").Value = "72:26:53"
$wsR4.Range("G5").Value = "70:04:26"

# --- Sheet R5: refresh elapsed-duration figure ---
$wsR5 = $wb.Worksheets.Item("R5")
$wsR5.Range("G2").Value = "431:01:12"

# --- Sheet R6: refresh elapsed-duration figure ---
$wsR6 = $wb.Worksheets.Item("R6")
$wsR6.Range("G2").Value = "71:33:30"
